# Edit script: shift weekly price records down by one row (new week inserted
# at the top of the data block, row 131), and append the previously-last
# record (old row 187) as a brand new row 188.
#
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F CategoriaID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo,
# L Precio maximo, M Precio promedio ponderado, N Unidad, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion.
#
# Only D (Fecha), J (Volumen), K, L, M (prices) and P (Precio $/Kg) actually
# vary row-to-row in the existing data block, so we snapshot those six
# columns for rows 131-187 first (before making any changes), then shift
# them down by one row. Row 131 receives the genuinely new record's date
# and volume. Row 188 is a brand-new row that duplicates the old row 187
# in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 131
$lastRow = 187
$newLastRow = 188

# Snapshot the full old rows 131..187 (columns A..R) before touching anything.
# $snapshot is 1-based: $snapshot[i, c] is row (firstRow + i - 1), column c.
$snapshot = $ws.Range("A$($firstRow):R$($lastRow)").Value2()

# Shift rows 187 down to 132 (process bottom-up isn't strictly required since
# we rely purely on the snapshot, but we still do it top-down for clarity).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $srcOffset = ($r - 1) - $firstRow + 1   # 1-based offset into $snapshot for row (r-1)
    $ws.Cells.Item($r, 4).Value2 = $snapshot[$srcOffset, 4]    # D Fecha
    $ws.Cells.Item($r, 10).Value2 = $snapshot[$srcOffset, 10]  # J Volumen
    $ws.Cells.Item($r, 11).Value2 = $snapshot[$srcOffset, 11]  # K Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $snapshot[$srcOffset, 12]  # L Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $snapshot[$srcOffset, 13]  # M Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value2 = $snapshot[$srcOffset, 16]  # P Precio $/Kg
}

# Row 131 gets the genuinely new record's date and volume.
$ws.Cells.Item($firstRow, 4).Value2 = 44489
$ws.Cells.Item($firstRow, 10).Value2 = 400

# New row 188 duplicates old row 187 (captured in snapshot's last row) in full.
$lastSnapshotOffset = $lastRow - $firstRow + 1
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newLastRow, $c).Value2 = $snapshot[$lastSnapshotOffset, $c]
}

# Preserve the date number format on the new row's Fecha cell (column D),
# matching the style used throughout the Fecha column.
$ws.Cells.Item($newLastRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
